$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")
$ws.Rows.Item(18).Insert()
Write-Host "Inserted row"
